# Add season-record columns (Wins, Losses, Ties) as AD, AE, AF
# next to the existing team/player stats table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Copy the formatting of the last existing header cell (AC1, style index 1:
# bold, bordered, centered) onto the three new header cells so they match
# the rest of the header row, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-52) ---
# Every player/team row gets the same season record: 84 wins, 78 losses, 0 ties.
$wins = 84
$losses = 78
$ties = 0

for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD = 30
    $ws.Cells.Item($r, 31).Value = $losses  # column AE = 31
    $ws.Cells.Item($r, 32).Value = $ties    # column AF = 32
}
